$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so the old CSE/B.Tech faculty-course rows are gone.
$ws.Cells.Clear()

# ---- Header row (faculty contact sheet instead of course list) ----
$ws.Range("A1").Value = "faculty_name"
$ws.Range("B1").Value = "email"
$ws.Range("C1").Value = "phone"
$ws.Range("D1").Value = "discipline_id"
$ws.Range("E1").Value = "branch_id"

# Header A1 gets a centered / wrapped style
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108
$ws.Range("A1").WrapText = $true

# ---- Data row ----
$ws.Range("B2").Value = "exceel@gmail.com"
$ws.Range("A2").Value = "ME specialist"
$ws.Range("C2").Value = 7418529633
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 3

# Turn the email into a real mailto hyperlink (auto-applies the Hyperlink style)
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:exceel@gmail.com")

# The rest of the (now mostly empty) student rows keep the Hyperlink-style look in column B
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("B6").Style = "Hyperlink"

# ---- Column widths ----
$ws.Columns("A").ColumnWidth = 13.665
$ws.Columns("B").ColumnWidth = 19.33
$ws.Columns("C").ColumnWidth = 10.165
$ws.Columns("D").ColumnWidth = 14.33

# ---- Selection moves to F6 ----
$ws.Range("F6").Select()

Write-Output "done"
